$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(66.13536318407947,51.212,79.91999999999996),
    @(7.00907462686566,3.5,11.628),
    @(15.10718407960198,8.151999999999994,23.18),
    @(8.425830845771129,4.112000000000001,12.86799999999999),
    @(33.52193034825862,21.5,46.956),
    @(61.57703482587047,45.13199999999997,77.52000000000005),
    @(91.00686567164175,83.13999999999999,96.13599999999995),
    @(31.34748258706461,18.06,43.76399999999999),
    @(90.95800995024865,82.67600000000004,96.16800000000001),
    @(78.03156218905455,62.47600000000001,88.53999999999998),
    @(53.24171144278596,35.47999999999998,69.78799999999997),
    @(89.05251741293527,79.33599999999997,95.21600000000001),
    @(84.72591044776105,73.49199999999995,92.83199999999999),
    @(9.366427860696511,4.883999999999999,15.04000000000001),
    @(33.41556218905465,21.12799999999999,46.768),
    @(52.99205970149245,36.45999999999999,69.54400000000001),
    @(15.06483582089552,8.115999999999993,23.068),
    @(78.82869651741282,63.696,89.236),
    @(36.61416915422874,22.72,52.78),
    @(20.0030646766169,11.43199999999999,29.38400000000001),
    @(91.35667661691535,82.18000000000004,96.47199999999997),
    @(10.43315422885571,5.752000000000002,16.76400000000001),
    @(33.45353233830838,21.164,46.756),
    @(29.02539303482581,17.56799999999999,42.42800000000001),
    @(19.64055721393035,12.012,28.93599999999998),
    @(89.19104477611931,77.73200000000006,95.47199999999999),
    @(64.38459701492525,47.048,78.25199999999994),
    @(10.62479601990049,5.872000000000003,16.768),
    @(51.18905472636802,34.492,66.67599999999995),
    @(84.55832835820891,72.86800000000005,92.38),
    @(88.19299502487554,77.64000000000001,94.36800000000001),
    @(84.90750248756206,73.56799999999996,92.72800000000007),
    @(69.08981094527357,52.488,82.62799999999994),
    @(34.9180895522387,21.54000000000001,48.288),
    @(61.06722388059693,44.712,76.18400000000004),
    @(21.44835820895515,12.49199999999999,33.02799999999999),
    @(85.28774129353222,74.54000000000002,93.12400000000001),
    @(55.25084577114425,39.33599999999997,70.78000000000003),
    @(17.57880597014923,10.46000000000001,26.48399999999998),
    @(58.81838805970135,42.38400000000001,73.32400000000001),
    @(78.35008955223866,63.82399999999998,88.81200000000003),
    @(51.06260696517399,34.80400000000002,66.41199999999998),
    @(72.37426865671637,54.69600000000002,84.66000000000005),
    @(79.97120398009936,68.38799999999996,89.65599999999996),
    @(83.14489552238786,70.93999999999998,91.79600000000001),
    @(91.32306467661674,82.69999999999996,96.34800000000003),
    @(53.10093532338296,35.29199999999999,69.27199999999995),
    @(16.03295522388058,9.396000000000006,24.62),
    @(21.07265671641781,11.736,32.42799999999998),
    @(6.999820895522375,3.452000000000001,11.656),
    @(81.76724378109454,66.90399999999998,90.46799999999998),
    @(55.13892537313427,39.12799999999998,71.21600000000001),
    @(9.366447761194022,4.895999999999999,15.03200000000001),
    @(19.65818905472637,11.54,29.65599999999998),
    @(84.38185074626848,71.59600000000002,92.07200000000006),
    @(88.13755223880592,77.77200000000001,94.46400000000003),
    @(21.55313432835809,12.688,33.76799999999999),
    @(76.3960995024874,60.236,87.78400000000002),
    @(10.40179104477611,5.344000000000001,16.33200000000001),
    @(85.19186069651734,73.83199999999999,92.71200000000003),
    @(78.01452736318393,63.32,88.73999999999998),
    @(85.68069651741276,75.07199999999995,93.15199999999993),
    @(10.53799004975124,5.704000000000001,16.492),
    @(87.42234825870635,77.71600000000002,93.82800000000002),
    @(68.37795024875614,51.10400000000001,82.24399999999994),
    @(66.86718407960193,50.78,80.86),
    @(19.01323383084576,11.104,28.94799999999999),
    @(91.54979104477596,83.20400000000004,96.56399999999999),
    @(17.63389054726368,9.912000000000006,26.468),
    @(42.62117412935316,27.52000000000001,56.91200000000004),
    @(48.58551243781078,30.70800000000001,62.57200000000003)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
